$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '62.678.79'
$ws.Range("E2").Value = '  +1.02%  '
Set-TextValue $ws.Range("D3") '2.437.42'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue $ws.Range("D5") '567.22'
$ws.Range("E5").Value = '  +0.62%  '
Set-TextValue $ws.Range("D6") '145.37'
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("E7").Value = '  -0.01%  '
Set-TextValue $ws.Range("D8") '0.534'
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("E10").Value = '  +0.48%  '
Set-TextValue $ws.Range("D11") '5.25'
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("E12").Value = '  +0.99%  '
Set-TextValue $ws.Range("D13") '26.83'
$ws.Range("E13").Value = '  +4.61%  '
$ws.Range("E14").Value = '  +5.47%  '
$ws.Range("E15").Value = '  +1.10%  '
Set-TextValue $ws.Range("D16") '62.609.11'
$ws.Range("E16").Value = '  +0.91%  '
Set-TextValue $ws.Range("D17") '2.434.27'
$ws.Range("E17").Value = '  +1.06%  '
Set-TextValue $ws.Range("D18") '11.23'
$ws.Range("E18").Value = '  -0.62%  '
Set-TextValue $ws.Range("D19") '6.92'
$ws.Range("E19").Value = '  +0.76%  '
Set-TextValue $ws.Range("D20") '323.35'
$ws.Range("E20").Value = '  +0.44%  '
Set-TextValue $ws.Range("D21") '4.17'
$ws.Range("E21").Value = '  +0.79%  '
Set-TextValue $ws.Range("D22") '0.999'
$ws.Range("E22").Value = '  -0.15%  '
Set-TextValue $ws.Range("D23") '67.28'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("E24").Value = '  +2.57%  '
Set-TextValue $ws.Range("D25") '8.70'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  +7.89%  '
Set-TextValue $ws.Range("D27") '566.11'
$ws.Range("E27").Value = '  -1.08%  '
Set-TextValue $ws.Range("D28") '2.553.63'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("E29").Value = '  -0.28%  '
Set-TextValue $ws.Range("D30") '8.40'
$ws.Range("E30").Value = '  +2.64%  '
Set-TextValue $ws.Range("D31") '1.45'
$ws.Range("E31").Value = '  +2.76%  '
Set-TextValue $ws.Range("D32") '0.147'
$ws.Range("E32").Value = '  -0.45%  '
Set-TextValue $ws.Range("D33") '1.88'
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("E34").Value = '  +1.53%  '
Set-TextValue $ws.Range("D35") '4.86'
$ws.Range("E35").Value = '  +3.71%  '
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("E37").Value = '  +0.78%  '
Set-TextValue $ws.Range("D38") '5.45'
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("E39").Value = '  +0.94%  '
Set-TextValue $ws.Range("D40") '148.27'
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("E41").Value = '  +2.09%  '
$ws.Range("E42").Value = '  +0.42%  '
Set-TextValue $ws.Range("D43") '2.42'
$ws.Range("E43").Value = '  +6.08%  '
Set-TextValue $ws.Range("D44") '148.44'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("E46").Value = '  +1.05%  '
Set-TextValue $ws.Range("D47") '20.50'
$ws.Range("E47").Value = '  +2.78%  '
Set-TextValue $ws.Range("D48") '0.599'
$ws.Range("E48").Value = '  +1.43%  '
Set-TextValue $ws.Range("D49") '0.0231'
$ws.Range("E49").Value = '  +2.85%  '
Set-TextValue $ws.Range("D50") '0.0927'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("E51").Value = '  +1.08%  '
